$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New section: "European Union" (rows 17-19), mirroring the existing
# "France" section (rows 12-14) that precedes it.
# ---------------------------------------------------------------------------

# Row 17: section header ("European Union" merged across C17:E17, with the
# "Rate necessary" label in G17) - styled the same way as row 12.
$ws.Range("C17").Value = "European Union"
$ws.Range("C17:E17").Merge()
$ws.Range("C17:E17").HorizontalAlignment = -4108
$ws.Range("F17").NumberFormat = "0"
$ws.Range("G17").Value = "Rate necessary"

# Row 18: data row for the European Union.
$ws.Range("A18").Value = "p97p100"
$ws.Range("B18").Value = 2022
$ws.Range("C18").Value = 119588.4
$ws.Range("C18").NumberFormat = "0"
$ws.Range("D18").Formula = "=C18/12"
$ws.Range("D18").NumberFormat = "0"
$ws.Range("E18").Value = 240624.5
$ws.Range("E18").NumberFormat = "0"
$ws.Range("F18:F19").Formula = "=E18/12"
$ws.Range("F18:F19").NumberFormat = "0"
$ws.Range("G18").Formula = "=(100/3)*E19/(E18-C18)"
$ws.Range("G18").NumberFormat = "0.0%"

# Row 19: "To be raised:" label (same text/shared string as D14) plus the
# supporting value used by G18.
$ws.Range("D19").Value = "To be raised:"
$ws.Range("E19").Value = 296.95
$ws.Range("E19").NumberFormat = "0"

# Move the selection to match the author's last position.
$ws.Range("G19").Select()
